$d = $word.ActiveDocument

# The document's sole paragraph holds two leftover placeholder text runs
# ("sd" and "fsd", both 8pt) that sandwich the _GoBack bookmark pair:
#   [sd][bookmarkStart][bookmarkEnd][fsd]
# These were stray/garbled text left in the template and must be cleared
# out so the space-group value can later be inserted as a properly
# formatted (MathML) field, while keeping the bookmark itself intact.
$para = $d.Paragraphs(1)
$paraStart = $para.Range.Start
# Range.End includes the trailing paragraph mark; back up one character to
# land on the end of the actual run text ("...fsd").
$textEnd = $para.Range.End - 1

# Delete the trailing run ("fsd") first so the leading run's offsets stay
# valid, then delete the leading run ("sd").
$d.Range($textEnd - 3, $textEnd).Delete()
$d.Range($paraStart, $paraStart + 2).Delete()
